$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text (volume/number, week date range)
$ws.Range("A8").Value = 'Volume 32   Number  42'
$ws.Range("C9").Value = 'Report Covering the Week  10/13/2025  Through  10/19/2025'

# Cells whose displayed TYPE changes (number <-> text placeholder) - fix value then restore exact number format via paste-format from a stable donor cell
$ws.Range("C14").Value = 2
$ws.Range("D15").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = '''0'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = '***.*'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = 4
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = -50
$ws.Range("H14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = '''0'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = '***.*'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = '''0'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = '***.*'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("D31").Value = '''0'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = '***.*'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("G33").Value = '''0'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("G33").PasteSpecial(-4122) | Out-Null
$ws.Range("H33").Value = '***.*'
$ws.Range("C31").Copy() | Out-Null
$ws.Range("H33").PasteSpecial(-4122) | Out-Null

# Remaining numeric/text value-only updates (format/style unchanged)
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = -20
$ws.Range("I14").Value = 30
$ws.Range("K14").Value = -3.225806451612
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = -49.152542372881
$ws.Range("N14").Value = -75.409836065573
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 100
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = -5.263157894736
$ws.Range("I15").Value = 183
$ws.Range("J15").Value = 156
$ws.Range("K15").Value = 17.307692307692
$ws.Range("L15").Value = 55.084745762711
$ws.Range("M15").Value = 71.028037383177
$ws.Range("N15").Value = -35.563380281690
$ws.Range("C16").Value = 26
$ws.Range("D16").Value = 28
$ws.Range("E16").Value = -7.142857142857
$ws.Range("F16").Value = 112
$ws.Range("G16").Value = 117
$ws.Range("H16").Value = -4.273504273504
$ws.Range("I16").Value = 987
$ws.Range("J16").Value = 1140
$ws.Range("K16").Value = -13.421052631578
$ws.Range("L16").Value = -9.615384615384
$ws.Range("M16").Value = -41.804245283018
$ws.Range("N16").Value = -84.124175647418
$ws.Range("C17").Value = 52
$ws.Range("D17").Value = 74
$ws.Range("E17").Value = -29.729729729729
$ws.Range("F17").Value = 228
$ws.Range("G17").Value = 253
$ws.Range("H17").Value = -9.881422924901
$ws.Range("I17").Value = 2416
$ws.Range("J17").Value = 2574
$ws.Range("K17").Value = -6.138306138306
$ws.Range("L17").Value = 5.134899912967
$ws.Range("M17").Value = 77.908689248895
$ws.Range("N17").Value = -20.447810339150
$ws.Range("C18").Value = 18
$ws.Range("D18").Value = 29
$ws.Range("E18").Value = -37.931034482758
$ws.Range("F18").Value = 59
$ws.Range("G18").Value = 95
$ws.Range("H18").Value = -37.894736842105
$ws.Range("I18").Value = 819
$ws.Range("J18").Value = 869
$ws.Range("K18").Value = -5.753739930955
$ws.Range("L18").Value = -10.491803278688
$ws.Range("M18").Value = -51.162790697674
$ws.Range("N18").Value = -88.826739427012
$ws.Range("C19").Value = 75
$ws.Range("D19").Value = 68
$ws.Range("E19").Value = 10.294117647058
$ws.Range("F19").Value = 310
$ws.Range("G19").Value = 258
$ws.Range("H19").Value = 20.155038759689
$ws.Range("I19").Value = 2749
$ws.Range("J19").Value = 2667
$ws.Range("K19").Value = 3.074615673040
$ws.Range("L19").Value = -6.908228919742
$ws.Range("M19").Value = 21.101321585903
$ws.Range("N19").Value = -58.499396135265
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 57
$ws.Range("E20").Value = -19.298245614035
$ws.Range("F20").Value = 162
$ws.Range("G20").Value = 201
$ws.Range("H20").Value = -19.402985074626
$ws.Range("I20").Value = 1548
$ws.Range("J20").Value = 1630
$ws.Range("K20").Value = -5.030674846625
$ws.Range("L20").Value = 5.091649694501
$ws.Range("M20").Value = 8.860759493670
$ws.Range("N20").Value = -89.920562573251
$ws.Range("C21").Value = 223
$ws.Range("D21").Value = 258
$ws.Range("E21").Value = -13.565891472868
$ws.Range("F21").Value = 893
$ws.Range("G21").Value = 948
$ws.Range("H21").Value = -5.801687763713
$ws.Range("I21").Value = 8732
$ws.Range("J21").Value = 9067
$ws.Range("K21").Value = -3.694717105988
$ws.Range("L21").Value = -1.489169675090
$ws.Range("M21").Value = 1.664920246827
$ws.Range("N21").Value = -77.594170173457
$ws.Range("F22").Value = 9
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 80
$ws.Range("I22").Value = 92
$ws.Range("J22").Value = 106
$ws.Range("K22").Value = -13.207547169811
$ws.Range("L22").Value = -4.166666666666
$ws.Range("M22").Value = -6.122448979591
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 16
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 179
$ws.Range("J23").Value = 182
$ws.Range("K23").Value = -1.648351648351
$ws.Range("L23").Value = -5.789473684210
$ws.Range("M23").Value = 47.933884297520
$ws.Range("C24").Value = 194
$ws.Range("D24").Value = 182
$ws.Range("E24").Value = 6.593406593406
$ws.Range("F24").Value = 776
$ws.Range("G24").Value = 802
$ws.Range("H24").Value = -3.241895261845
$ws.Range("I24").Value = 7457
$ws.Range("J24").Value = 7540
$ws.Range("K24").Value = -1.100795755968
$ws.Range("L24").Value = 0.743042420967
$ws.Range("M24").Value = 50.070436707587
$ws.Range("C25").Value = 70
$ws.Range("D25").Value = 93
$ws.Range("E25").Value = -24.731182795698
$ws.Range("F25").Value = 308
$ws.Range("G25").Value = 403
$ws.Range("H25").Value = -23.573200992555
$ws.Range("I25").Value = 2973
$ws.Range("J25").Value = 3307
$ws.Range("K25").Value = -10.099788327789
$ws.Range("L25").Value = 15.366705471478
$ws.Range("C26").Value = 82
$ws.Range("D26").Value = 84
$ws.Range("E26").Value = -2.380952380952
$ws.Range("F26").Value = 393
$ws.Range("G26").Value = 412
$ws.Range("H26").Value = -4.611650485436
$ws.Range("I26").Value = 3979
$ws.Range("J26").Value = 4174
$ws.Range("K26").Value = -4.671777671298
$ws.Range("L26").Value = 8.301578660860
$ws.Range("M26").Value = 8.419618528610
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 24
$ws.Range("H27").Value = 4.166666666666
$ws.Range("I27").Value = 229
$ws.Range("J27").Value = 245
$ws.Range("K27").Value = -6.530612244897
$ws.Range("L27").Value = 10.096153846153
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = -16.666666666666
$ws.Range("F28").Value = 36
$ws.Range("G28").Value = 39
$ws.Range("H28").Value = -7.692307692307
$ws.Range("I28").Value = 363
$ws.Range("J28").Value = 359
$ws.Range("K28").Value = 1.114206128133
$ws.Range("L28").Value = 12.732919254658
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = -54.545454545454
$ws.Range("I29").Value = 59
$ws.Range("K29").Value = -37.234042553191
$ws.Range("L29").Value = -20.270270270270
$ws.Range("M29").Value = -62.420382165605
$ws.Range("N29").Value = -85.783132530120
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 9
$ws.Range("H30").Value = -44.444444444444
$ws.Range("I30").Value = 47
$ws.Range("K30").Value = -38.157894736842
$ws.Range("L30").Value = -16.071428571428
$ws.Range("M30").Value = -62.4
$ws.Range("N30").Value = -87.5
$ws.Range("I31").Value = 32
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = -31.914893617021
$ws.Range("F33").Value = 2
$ws.Range("L33").Value = -16.666666666666

$excel.CutCopyMode = 0
